# Auto-generated Excel COM-interop script to refresh market-price derived columns (H-N)
# across all 8 sheets, per commit "chore: update Sheets via scheduled runner".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 116.2
$ws.Range("I4").Value = 84.666664
$ws.Range("K4").Value = 84.666664
$ws.Range("M4").Value = 29.333336
$ws.Range("H9").Value = 206.44444
$ws.Range("I9").Value = 251.83333
$ws.Range("J9").Value = 115.666664
$ws.Range("K9").Value = 251.83333
$ws.Range("L9").Value = 115.666664
$ws.Range("M9").Value = -82.83332999999999
$ws.Range("N9").Value = -453.666664
$ws.Range("H17").Value = 1497.3
$ws.Range("J17").Value = 1497.1578
$ws.Range("L17").Value = 4491.4734
$ws.Range("N17").Value = -4827.4734
$ws.Range("H20").Value = 34085.168
$ws.Range("I20").Value = 26127.75
$ws.Range("K20").Value = 26127.75
$ws.Range("M20").Value = -25897.75
$ws.Range("H35").Value = 34085.168
$ws.Range("I35").Value = 26127.75
$ws.Range("K35").Value = 26127.75
$ws.Range("M35").Value = -25748.75
$ws.Range("H96").Value = 2333.1035
$ws.Range("I96").Value = 2776.8
$ws.Range("J96").Value = 1347.1111
$ws.Range("K96").Value = 8330.400000000001
$ws.Range("L96").Value = 4041.3333
$ws.Range("M96").Value = -6957.400000000001
$ws.Range("N96").Value = -6787.3333
$ws.Range("H98").Value = 1214.7778
$ws.Range("I98").Value = 1580
$ws.Range("J98").Value = 888
$ws.Range("K98").Value = 1580
$ws.Range("L98").Value = 888
$ws.Range("M98").Value = -82
$ws.Range("N98").Value = -3884
$ws.Range("H122").Value = 1214.7778
$ws.Range("I122").Value = 1580
$ws.Range("J122").Value = 888
$ws.Range("K122").Value = 4740
$ws.Range("L122").Value = 2664
$ws.Range("M122").Value = -2290
$ws.Range("N122").Value = -7564

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 31094726
$ws.Range("I2").Value = 44049280
$ws.Range("J2").Value = 3799.6
$ws.Range("K2").Value = 44049280
$ws.Range("L2").Value = 3799.6
$ws.Range("M2").Value = -44049167
$ws.Range("N2").Value = -4025.6
$ws.Range("H32").Value = 7078.9272
$ws.Range("I32").Value = 4415.5713
$ws.Range("J32").Value = 28829.666
$ws.Range("K32").Value = 4415.5713
$ws.Range("L32").Value = 28829.666
$ws.Range("M32").Value = -4128.5713
$ws.Range("N32").Value = -29403.666
$ws.Range("H61").Value = 4448.074
$ws.Range("I61").Value = 3732.75
$ws.Range("J61").Value = 6491.857
$ws.Range("K61").Value = 3732.75
$ws.Range("L61").Value = 6491.857
$ws.Range("M61").Value = -3520.75
$ws.Range("N61").Value = -6915.857
$ws.Range("H63").Value = 3793.5
$ws.Range("I63").Value = 3632.1667
$ws.Range("K63").Value = 3632.1667
$ws.Range("M63").Value = -2946.1667
$ws.Range("H66").Value = 3793.5
$ws.Range("I66").Value = 3632.1667
$ws.Range("K66").Value = 18160.8335
$ws.Range("M66").Value = -14728.8335
$ws.Range("H97").Value = 35751040
$ws.Range("I97").Value = 47620360
$ws.Range("J97").Value = 143081.14
$ws.Range("K97").Value = 47620360
$ws.Range("L97").Value = 143081.14
$ws.Range("M97").Value = -47619864
$ws.Range("N97").Value = -144073.14
$ws.Range("H116").Value = 31094726
$ws.Range("I116").Value = 44049280
$ws.Range("J116").Value = 3799.6
$ws.Range("K116").Value = 44049280
$ws.Range("L116").Value = 3799.6
$ws.Range("M116").Value = -44046986
$ws.Range("N116").Value = -8387.6
$ws.Range("H132").Value = 937.5714
$ws.Range("I132").Value = 710.5833
$ws.Range("K132").Value = 2131.7499
$ws.Range("M132").Value = 398.2501000000002
$ws.Range("H135").Value = 109593
$ws.Range("J135").Value = 109593
$ws.Range("L135").Value = 109593
$ws.Range("N135").Value = -119733
$ws.Range("H136").Value = 4448.074
$ws.Range("I136").Value = 3732.75
$ws.Range("J136").Value = 6491.857
$ws.Range("K136").Value = 11198.25
$ws.Range("L136").Value = 19475.571
$ws.Range("M136").Value = -8648.25
$ws.Range("N136").Value = -24575.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 31094726
$ws.Range("I3").Value = 44049280
$ws.Range("J3").Value = 3799.6
$ws.Range("K3").Value = 44049280
$ws.Range("L3").Value = 3799.6
$ws.Range("M3").Value = -44049166
$ws.Range("N3").Value = -4027.6
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""
$ws.Range("H86").Value = 1961.375
$ws.Range("I86").Value = 1129.1428
$ws.Range("K86").Value = 1129.1428
$ws.Range("M86").Value = -6.142800000000079
$ws.Range("H89").Value = 1961.375
$ws.Range("I89").Value = 1129.1428
$ws.Range("K89").Value = 5645.714
$ws.Range("M89").Value = -29.71399999999994
$ws.Range("H94").Value = 933.7059
$ws.Range("I94").Value = 750.3077
$ws.Range("J94").Value = 1529.75
$ws.Range("K94").Value = 750.3077
$ws.Range("L94").Value = 1529.75
$ws.Range("M94").Value = -299.3077
$ws.Range("N94").Value = -2431.75
$ws.Range("H107").Value = 3849.8235
$ws.Range("I107").Value = 3243
$ws.Range("K107").Value = 3243
$ws.Range("M107").Value = -1323

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3336301.2
$ws.Range("I86").Value = 4764638
$ws.Range("K86").Value = 4764638
$ws.Range("M86").Value = -4763515
$ws.Range("H89").Value = 3336301.2
$ws.Range("I89").Value = 4764638
$ws.Range("K89").Value = 23823190
$ws.Range("M89").Value = -23817574
$ws.Range("H107").Value = 1488.619
$ws.Range("I107").Value = 1444.4667
$ws.Range("K107").Value = 1444.4667
$ws.Range("M107").Value = 475.5333000000001
$ws.Range("H108").Value = 26617.7
$ws.Range("I108").Value = 13356.286
$ws.Range("J108").Value = 57561
$ws.Range("K108").Value = 13356.286
$ws.Range("L108").Value = 57561
$ws.Range("M108").Value = -9516.286
$ws.Range("N108").Value = -65241
$ws.Range("H132").Value = 6693.4707
$ws.Range("I132").Value = 4149.1665
$ws.Range("K132").Value = 12447.4995
$ws.Range("M132").Value = -9917.499500000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 405.57144
$ws.Range("J103").Value = 371
$ws.Range("L103").Value = 1113
$ws.Range("N103").Value = -2871
$ws.Range("H122").Value = 1693.4375
$ws.Range("J122").Value = 1739.7333
$ws.Range("L122").Value = 15657.5997
$ws.Range("N122").Value = -20557.5997
$ws.Range("H131").Value = 563174.5
$ws.Range("I131").Value = 4246555
$ws.Range("K131").Value = 12739665
$ws.Range("M131").Value = -12734625
$ws.Range("H132").Value = 2362.3667
$ws.Range("J132").Value = 2933.65
$ws.Range("L132").Value = 26402.85
$ws.Range("N132").Value = -31462.85

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1203.2632
$ws.Range("I107").Value = 1635.75
$ws.Range("K107").Value = 1635.75
$ws.Range("M107").Value = 284.25
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""
$ws.Range("H132").Value = 2736.3845
$ws.Range("I132").Value = 1278.8334
$ws.Range("K132").Value = 3836.5002
$ws.Range("M132").Value = -1306.5002
$ws.Range("H136").Value = 96500
$ws.Range("J136").Value = 96500
$ws.Range("L136").Value = 289500
$ws.Range("N136").Value = -294600
$ws.Range("H137").Value = 103593
$ws.Range("J137").Value = 103593
$ws.Range("L137").Value = 103593
$ws.Range("N137").Value = -113793
$ws.Range("H138").Value = 145000
$ws.Range("J138").Value = 145000
$ws.Range("L138").Value = 145000
$ws.Range("N138").Value = -155280
$ws.Range("H141").Value = 42398.4
$ws.Range("J141").Value = 42398.4
$ws.Range("L141").Value = 42398.4
$ws.Range("N141").Value = -52758.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5000831.5
$ws.Range("I16").Value = 7813115
$ws.Range("K16").Value = 7813115
$ws.Range("M16").Value = -7812945
$ws.Range("H22").Value = 846.2
$ws.Range("I22").Value = 859.8
$ws.Range("J22").Value = 819
$ws.Range("K22").Value = 859.8
$ws.Range("L22").Value = 819
$ws.Range("M22").Value = -564.8
$ws.Range("N22").Value = -1409
$ws.Range("H27").Value = 846.2
$ws.Range("I27").Value = 859.8
$ws.Range("J27").Value = 819
$ws.Range("K27").Value = 859.8
$ws.Range("L27").Value = 819
$ws.Range("M27").Value = -752.8
$ws.Range("N27").Value = -1033
$ws.Range("H40").Value = 4962.8335
$ws.Range("I40").Value = 3444.25
$ws.Range("K40").Value = 3444.25
$ws.Range("M40").Value = -3308.25
$ws.Range("H68").Value = 37040260
$ws.Range("I68").Value = 47621620
$ws.Range("K68").Value = 47621620
$ws.Range("M68").Value = -47620871
$ws.Range("H71").Value = 37040260
$ws.Range("I71").Value = 47621620
$ws.Range("K71").Value = 238108100
$ws.Range("M71").Value = -238104356
$ws.Range("H122").Value = 4820
$ws.Range("I122").Value = 3548
$ws.Range("K122").Value = 10644
$ws.Range("M122").Value = -8194

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3563.9412
$ws.Range("I96").Value = 4010.9
$ws.Range("J96").Value = 2925.4285
$ws.Range("K96").Value = 4010.9
$ws.Range("L96").Value = 2925.4285
$ws.Range("M96").Value = -2637.9
$ws.Range("N96").Value = -5671.4285
$ws.Range("H107").Value = 1174.8572
$ws.Range("I107").Value = 1248.8182
$ws.Range("K107").Value = 3746.4546
$ws.Range("M107").Value = -1826.4546
$ws.Range("H141").Value = 120465
$ws.Range("J141").Value = 120465
$ws.Range("L141").Value = 120465
$ws.Range("N141").Value = -130825

Write-Output "Applied 252 cell updates across 8 sheets."